# PracticaGitHub.docx edit — "Modificación en main para rebase"
#
# 1. Insert a new paragraph right after the "merge" section's second
#    screenshot (the last inline picture in the document) containing the
#    heading text "4. Introducción a rebase".
# 2. Mark the two image-only paragraphs that precede that new heading
#    (the two screenshots under the "merge" section) as NoProofing, which
#    serializes as <w:rPr><w:noProof/></w:rPr> on their runs.
#
# (Order matters: NoProofing is inherited by InsertParagraphAfter, so the
#  new heading paragraph is created before the pictures are marked
#  NoProofing, keeping the new run's rPr empty as in the target diff.)

$d = $word.ActiveDocument

# The last inline picture in the document is the second "merge" screenshot
# (wp14:anchorId 5781732B); the new heading goes right after its paragraph.
$lastShape = $d.InlineShapes.Item($d.InlineShapes.Count)
$lastPara = $lastShape.Range.Paragraphs.Item(1)

$lastPara.Range.InsertParagraphAfter()
$insertedPara = $d.Paragraphs.Item($lastPara.Index + 1)
$insertedPara.Range.Text = "4. Introducción a rebase"

# The two "merge" screenshots (heights 861060 EMU and 1724660 EMU) get
# NoProofing on their paragraph's run, producing <w:noProof/> in rPr.
$targetHeightsEmu = @(861060, 1724660)
foreach ($shape in $d.InlineShapes) {
    $heightEmu = [Math]::Round($shape.Height * 12700)
    if ($targetHeightsEmu -contains $heightEmu) {
        $shape.Range.Paragraphs.Item(1).Range.NoProofing = 1
    }
}
